# "Refined metadata to be additional tab"
#
# 1. Update the panel_query_time-ish timestamps in the "data" sheet's
#    time_taken column (F2:F15) to the refreshed run's timestamps.
# 2. Add a new "metadata" worksheet (placed after "data") describing the
#    PanelApp query that produced the "data" sheet.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------
# 1. Refresh the per-row query timestamps on the "data" sheet.
# ---------------------------------------------------------------------
$timeTaken = @(
    "2021-10-05 14:20:20.408794",
    "2021-10-05 14:20:20.408801",
    "2021-10-05 14:20:20.408805",
    "2021-10-05 14:20:20.408807",
    "2021-10-05 14:20:20.408810",
    "2021-10-05 14:20:20.408813",
    "2021-10-05 14:20:20.408815",
    "2021-10-05 14:20:20.408818",
    "2021-10-05 14:20:20.408821",
    "2021-10-05 14:20:20.408823",
    "2021-10-05 14:20:20.408826",
    "2021-10-05 14:20:20.408828",
    "2021-10-05 14:20:20.408831",
    "2021-10-05 14:20:20.408833"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $data.Range("F$row").Value = $timeTaken[$i]
}

# ---------------------------------------------------------------------
# 2. Add the "metadata" sheet right after "data".
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Add([System.Type]::Missing, $data)
$meta.Name = "metadata"

# Match the look of the "data" sheet's header row (bold / bordered /
# centered) by copying its formatting onto the new header cells.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Match the "data" sheet's index-column style (A2) for metadata!A2.
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Familial prostate cancer"
$meta.Range("C2").Value = 318
$meta.Range("E2").Value = "2017-11-05T02:37:20.419988Z"
$meta.Range("F2").Value = "2021-10-05 14:20:20.405154"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/318/?format=json"

# "1.1" needs to land as literal text (not be auto-coerced to the number
# 1.1) while still leaving the cell in the sheet's default (unstyled)
# format -- stage it in a scratch cell formatted as Text, then bring over
# only the value with Paste Special so no number-format style sticks to
# D2 itself.
$scratch = $meta.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1.1"
$scratch.Copy()
$meta.Range("D2").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

$meta.Range("A1").Select() | Out-Null
